$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the two FR00 / 2040 / hydrogen demand values (elysis increase)
$ws.Range("E9").Value = 1600
$ws.Range("E17").Value = 8000

# Re-point the AutoFilter on the Year column (col 4) from 2030 to 2040.
# xlFilterValues = 7 -> emits <filters><filter val="..."/></filters>
# and also recomputes which rows are hidden/visible to match the new filter.
$ws.Range("A1:E31").AutoFilter(4, "2040", 7)

# Move the active cell / selection to E16 (matches the saved view state).
$ws.Range("E16").Select()
